# Updated Hybrid framework with all runmode set to N
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Set Runmode column (B2:B5) to "N" for all test cases
$ws.Range("B2:B5").Value = "N"

# Match the selection state recorded in the saved file
$ws.Range("B2:B5").Select()
